$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 602; existing rows 602..652 shift down to 603..653
$ws.Rows.Item(602).Insert()

# Populate the new row 602 with the new price-report record
$ws.Cells.Item(602, 1).Value = 5
$ws.Cells.Item(602, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(602, 3).Value = "Maule"
$ws.Cells.Item(602, 4).Value = 45008
$ws.Cells.Item(602, 5).Value = 7
$ws.Cells.Item(602, 6).Value = "Fruta"
$ws.Cells.Item(602, 7).Value = 100103
$ws.Cells.Item(602, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(602, 9).Value = 100103004
$ws.Cells.Item(602, 10).Value = "Durazno"
$ws.Cells.Item(602, 11).Value = "September Sweet"
$ws.Cells.Item(602, 12).Value = "Extra (doble especial)"
$ws.Cells.Item(602, 13).Value = 100
$ws.Cells.Item(602, 14).Value = 16000
$ws.Cells.Item(602, 15).Value = 16000
$ws.Cells.Item(602, 16).Value = 16000
$ws.Cells.Item(602, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(602, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(602, 19).Value = 1067
$ws.Cells.Item(602, 20).Value = 15
